# Applies the "also added some info to support diagram" edit:
#   1. Insert a new sentence right after "...the speaker system built by
#      team one." in the Design section.
#   2. Move the hidden "_GoBack" bookmark from its old location (end of the
#      "...and an 8Ohm speaker." paragraph, in the Results section) to the
#      end of the newly inserted sentence - this is what Word itself does
#      automatically when the last edit happens at a new location.

$d = $word.ActiveDocument

$newSentence = " The design consists of multiple modules. So that each of these modules could be tested individually and then if they all meet the specifications the active loudspeaker was assembled by connecting all the modules."

# --- Step 1: locate the insertion point and insert the new sentence -----
$rng = $d.Content
[void]$rng.Find.Execute("the speaker system built by team one.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Append a throw-away marker character after the real sentence. This lets
# us anchor the new "_GoBack" bookmark around a real character (inside the
# paragraph) instead of a zero-length range sitting exactly on the
# paragraph-mark boundary, which this host mis-resolves back to the start
# of the document. We delete the marker once the bookmark is in place,
# leaving a clean, truly zero-length bookmark right after the real text.
$rng.InsertAfter($newSentence + "X")

# Tag the newly inserted run with the same "en-US" language the rest of
# the paragraph already uses. (Every other run in this paragraph is
# already w:lang w:val="en-US", so this leaves them textually unchanged.)
$rng.LanguageID = "en-US"

# --- Step 2: remove the old _GoBack bookmark -----------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# --- Step 3: re-create _GoBack at the end of the new sentence -----------
$endPos = $rng.End
$markerRng = $d.Range($endPos - 1, $endPos)
$d.Bookmarks.Add("_GoBack", $markerRng)

# Remove the throw-away marker character now that the bookmark anchors
# correctly right after the real inserted text.
$markerRng2 = $d.Range($endPos - 1, $endPos)
$markerRng2.Text = ""
